$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.399.81'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '3.462.06'
$ws.Range("E3").Value = '  -1.59%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'583.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = "'177.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("D7").Value = "'0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.63%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '3.460.72'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").Value = "'6.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = "'0.419"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").Value = '4.062.91'
$ws.Range("E13").Value = '  -1.59%  '
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("D15").Value = "'30.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = '66.296.84'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '3.485.76'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").Value = "'5.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("D20").Value = "'13.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("D21").Value = "'371.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("E22").Value = '  -3.17%  '
$ws.Range("D23").Value = "'73.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("D25").Value = "'0.537"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.54%  '
$ws.Range("E26").Value = '  +3.79%  '
$ws.Range("D27").Value = "'10.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("E28").Value = '  +2.76%  '
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").Value = "'5.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("D32").Value = "'23.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.86%  '
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").Value = "'7.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.54%  '
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("D36").Value = "'1.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("D37").Value = "'161.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("D38").Value = "'0.886"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("D39").Value = "'27.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.02%  '
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("D41").Value = '2.807.03'
$ws.Range("E41").Value = '  +2.92%  '
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = "'2.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("D44").Value = "'6.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("D45").Value = "'0.0693"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.46%  '
$ws.Range("D46").Value = "'25.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").Value = "'344.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.00%  '
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("D51").Value = "'31.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.20%  '
